$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: section header
$ws.Range("A31").Value = "Load times - Opera, Windows 7, cache disabled"

# Row 32: column headers
$ws.Range("B32").Value = "Time to load"
$ws.Range("C32").Value = "Time to display 25 comments"

# Row 33: Development
$ws.Range("A33").Value = "Development"
$ws.Range("B33").Value = "2800 ms"
$ws.Range("C33").Value = "~7000ms"

# Row 34: Production (Vulcanized)
$ws.Range("A34").Value = "Production (Vulcanized)"
$ws.Range("B34").Value = "1100 ms"
$ws.Range("C34").Value = "2800 ms"

# Row 36: Google Chrome section header
$ws.Range("A36").Value = "Google Chrome"

# Row 37: Dev
$ws.Range("A37").Value = "Dev"
$ws.Range("B37").Value = "2200 ms"
$ws.Range("C37").Value = "4200 ms"

# Row 38: Vulcanized
$ws.Range("A38").Value = "Vulcanized"
$ws.Range("B38").Value = "1150 ms"
$ws.Range("C38").Value = "2900 ms"

# Row 40: Firefox section header
$ws.Range("A40").Value = "Firefox"

# Row 41: Dev
$ws.Range("A41").Value = "Dev"
$ws.Range("B41").Value = "7400 ms"
$ws.Range("C41").Value = 12000

# Row 42: Vulcanized
$ws.Range("A42").Value = "Vulcanized"
$ws.Range("B42").Value = "6050 ms"
$ws.Range("C42").Value = 10500

# Update selection / view to match the author's final state
$ws.Range("C42").Select()
$excel.ActiveWindow.ScrollRow = 11
